$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Section 1: Corresponding Author block
#   "Samuel T Harrold Email: samuel.harrold@yuimedi.com Contact: https://us.yuimedi.com/contact-us/"
# becomes two paragraphs:
#   "Samuel T Harrold"                   (unchanged FirstParagraph style)
#   "Email: samuel.harrold@yuimedi.com"  (new BodyText paragraph; the "Contact:" sentence is dropped)
# ---------------------------------------------------------------------------

# Anchor on the unique "Corresponding Author" heading and compute the index of
# the paragraph that follows it, so we never touch the wrong occurrence of
# similar text elsewhere in the document.
$headingRange1 = $d.Content
$headingRange1.Find.Execute("Corresponding Author", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$before1 = $d.Range(0, $headingRange1.Start)
$headingIndex1 = $before1.Paragraphs.Count + 1
$contentIndex1 = $headingIndex1 + 1

$targetPara1 = $d.Paragraphs.Item($contentIndex1)
$targetRange1 = $targetPara1.Range
$targetRange1.Find.Execute(
    "Samuel T Harrold Email: samuel.harrold@yuimedi.com Contact: https://us.yuimedi.com/contact-us/",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Samuel T Harrold^pEmail: samuel.harrold@yuimedi.com", 2)

# The text is now split across two paragraphs; apply BodyText style to the
# second (newly created) one that holds the "Email: ..." line.
$d.Paragraphs.Item($contentIndex1 + 1).Style = "BodyText"

# ---------------------------------------------------------------------------
# Section 2: Signature block
#   "Samuel T Harrold Yuimedi, Inc." becomes two BodyText paragraphs:
#   "Samuel T Harrold"
#   "Yuimedi, Inc."
# ---------------------------------------------------------------------------

# Anchor on the unique "Sincerely," paragraph that immediately precedes the
# signature line.
$sincerelyRange = $d.Content
$sincerelyRange.Find.Execute("Sincerely,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$before2 = $d.Range(0, $sincerelyRange.Start)
$sincerelyIndex = $before2.Paragraphs.Count + 1
$contentIndex2 = $sincerelyIndex + 1

$targetPara2 = $d.Paragraphs.Item($contentIndex2)
$targetRange2 = $targetPara2.Range
$targetRange2.Find.Execute(
    "Samuel T Harrold Yuimedi, Inc.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Samuel T Harrold^pYuimedi, Inc.", 2)

# Apply BodyText style to the newly created second paragraph ("Yuimedi, Inc.").
$d.Paragraphs.Item($contentIndex2 + 1).Style = "BodyText"
